# ---- Step 1: "总计" sheet -- add the 2022-Q4 summary row ----
# Column A is a plain 0-based row index (0,1,2,...) and never moves; only the
# B/C/D (label/count/value) content slides down a row to make room for the new
# 2022-Q4 figures at the top, so values are written directly rather than via a
# row insert (which would drag the A-index along and desync it).
$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# Append the brand-new row 4, carrying forward what used to be row 3 data,
# and give its index cell (A4) the same bold+border style as A2/A3.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2020-Q4"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.95

# Row 3 becomes what used to be row 2s data (2021-Q4)
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.03

# Row 2 becomes the brand-new 2022-Q4 data
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 14
$totalSheet.Range("D2").Value = 3.04

# ---- Step 2: add the new "2022-Q4" sheet right after "总计" ----
$srcStyleSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# ---- Step 3: header row ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- Step 4: force text storage (preserve leading zeros / exact decimal text)
#              for columns B-F (all rows) and G (rows 2-13; rows 14-15 stay numeric) ----
$rngTextBF = $newSheet.Range("B2:F15")
$rngTextBF.NumberFormat = "@"
$rngTextG = $newSheet.Range("G2:G13")
$rngTextG.NumberFormat = "@"

# ---- Step 5: data rows ----
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "720001"
$newSheet.Range("C2").Value = "财通价值动量混合"
$newSheet.Range("D2").Value = "36.34"
$newSheet.Range("E2").Value = "79.52"
$newSheet.Range("F2").Value = "3.93"
$newSheet.Range("G2").Value = "1.4282"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "011201"
$newSheet.Range("C3").Value = "财通优势行业轮动混合A"
$newSheet.Range("D3").Value = "8.38"
$newSheet.Range("E3").Value = "88.19"
$newSheet.Range("F3").Value = "5.62"
$newSheet.Range("G3").Value = "0.4710"
$newSheet.Range("H3").Value = 3

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "501085"
$newSheet.Range("C4").Value = "财通科创主题灵活配置混合（LOF）"
$newSheet.Range("D4").Value = "4.19"
$newSheet.Range("E4").Value = "89.51"
$newSheet.Range("F4").Value = "7.05"
$newSheet.Range("G4").Value = "0.2954"
$newSheet.Range("H4").Value = 1

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "519909"
$newSheet.Range("C5").Value = "华安安顺灵活配置混合A"
$newSheet.Range("D5").Value = "9.68"
$newSheet.Range("E5").Value = "80.72"
$newSheet.Range("F5").Value = "2.65"
$newSheet.Range("G5").Value = "0.2565"
$newSheet.Range("H5").Value = 8

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "008983"
$newSheet.Range("C6").Value = "财通科技创新混合A"
$newSheet.Range("D6").Value = "2.87"
$newSheet.Range("E6").Value = "87.06"
$newSheet.Range("F6").Value = "4.10"
$newSheet.Range("G6").Value = "0.1177"
$newSheet.Range("H6").Value = 8

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "501032"
$newSheet.Range("C7").Value = "财通福盛多策略混合（LOF）A"
$newSheet.Range("D7").Value = "2.16"
$newSheet.Range("E7").Value = "87.95"
$newSheet.Range("F7").Value = "5.30"
$newSheet.Range("G7").Value = "0.1145"
$newSheet.Range("H7").Value = 4

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "009062"
$newSheet.Range("C8").Value = "财通智慧成长混合A"
$newSheet.Range("D8").Value = "2.05"
$newSheet.Range("E8").Value = "86.49"
$newSheet.Range("F8").Value = "4.67"
$newSheet.Range("G8").Value = "0.0957"
$newSheet.Range("H8").Value = 8

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "009063"
$newSheet.Range("C9").Value = "财通智慧成长混合C"
$newSheet.Range("D9").Value = "1.74"
$newSheet.Range("E9").Value = "86.49"
$newSheet.Range("F9").Value = "4.67"
$newSheet.Range("G9").Value = "0.0813"
$newSheet.Range("H9").Value = 8

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "008984"
$newSheet.Range("C10").Value = "财通科技创新混合C"
$newSheet.Range("D10").Value = "1.79"
$newSheet.Range("E10").Value = "87.06"
$newSheet.Range("F10").Value = "4.10"
$newSheet.Range("G10").Value = "0.0734"
$newSheet.Range("H10").Value = 8

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "006503"
$newSheet.Range("C11").Value = "财通集成电路产业股票C"
$newSheet.Range("D11").Value = "0.70"
$newSheet.Range("E11").Value = "80.87"
$newSheet.Range("F11").Value = "6.33"
$newSheet.Range("G11").Value = "0.0443"
$newSheet.Range("H11").Value = 2

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "006502"
$newSheet.Range("C12").Value = "财通集成电路产业股票A"
$newSheet.Range("D12").Value = "0.66"
$newSheet.Range("E12").Value = "80.87"
$newSheet.Range("F12").Value = "6.33"
$newSheet.Range("G12").Value = "0.0418"
$newSheet.Range("H12").Value = 2

$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "011202"
$newSheet.Range("C13").Value = "财通优势行业轮动混合C"
$newSheet.Range("D13").Value = "0.34"
$newSheet.Range("E13").Value = "88.19"
$newSheet.Range("F13").Value = "5.62"
$newSheet.Range("G13").Value = "0.0191"
$newSheet.Range("H13").Value = 3

$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "014628"
$newSheet.Range("C14").Value = "财通福盛多策略混合（LOF）C"
$newSheet.Range("D14").Value = "0.00"
$newSheet.Range("E14").Value = "87.95"
$newSheet.Range("F14").Value = "5.30"
$newSheet.Range("G14").Value = 0
$newSheet.Range("H14").Value = 4

$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "016564"
$newSheet.Range("C15").Value = "华安安顺灵活配置混合C"
$newSheet.Range("D15").Value = "0.00"
$newSheet.Range("E15").Value = "80.72"
$newSheet.Range("F15").Value = "2.65"
$newSheet.Range("G15").Value = 0
$newSheet.Range("H15").Value = 8

# ---- Step 6: drop the helper text-number-format now that values are typed ----
$rngTextBF.ClearFormats()
$rngTextG.ClearFormats()

# ---- Step 7: apply header / index-column styling (bold + border, cellXf "2")
#              from the existing quarterly sheet, matching its siblings ----
$srcStyleSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$srcStyleSheet.Range("A2").Copy()
$newSheet.Range("A2:A15").PasteSpecial(-4122)
